$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2' = '59.157.49'
    'E2' = '  -2.94%  '
    'D3' = '2.657.17'
    'E3' = '  -1.49%  '
    'E4' = '  -0.15%  '
    'D5' = '''523.75'
    'E5' = '  +0.24%  '
    'D6' = '''144.82'
    'E6' = '  -2.27%  '
    'E7' = '  +0.25%  '
    'E8' = '  -1.30%  '
    'D9' = '''7.04'
    'E9' = '  +9.40%  '
    'E10' = '  -3.74%  '
    'E11' = '  -2.25%  '
    'E12' = '  +1.74%  '
    'D13' = '3.119.77'
    'E13' = '  -1.72%  '
    'D14' = '59.173.15'
    'E14' = '  -3.02%  '
    'D15' = '''21.20'
    'E15' = '  -1.50%  '
    'E16' = '  -2.34%  '
    'D17' = '2.666.25'
    'E17' = '  -6.60%  '
    'D18' = '''340.58'
    'E18' = '  -4.13%  '
    'D19' = '''4.40'
    'E19' = '  -4.28%  '
    'E20' = '  -1.95%  '
    'D21' = '''6.37'
    'E21' = '  +0.00%  '
    'E22' = '  -0.06%  '
    'D23' = '''64.15'
    'E23' = '  +1.94%  '
    'E24' = '  -2.57%  '
    'E25' = '  -1.92%  '
    'D26' = '''0.999'
    'E26' = '  +0.85%  '
    'E27' = '  -2.78%  '
    'E28' = '  -2.54%  '
    'E29' = '  -1.57%  '
    'E30' = '  +0.08%  '
    'E31' = '  -0.28%  '
    'E32' = '  -2.00%  '
    'D33' = '''149.06'
    'E33' = '  -0.82%  '
    'D34' = '''4.17'
    'E34' = '  -1.11%  '
    'E35' = '  -2.71%  '
    'E36' = '  -5.18%  '
    'D37' = '''0.885'
    'E37' = '  +0.51%  '
    'B38' = 'Stacks'
    'C38' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'D38' = '''1.49'
    'E38' = '  -5.82%  '
    'B39' = 'OKB'
    'C39' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D39' = '''36.69'
    'E39' = '  -0.37%  '
    'E40' = '  -3.74%  '
    'E41' = '  +0.31%  '
    'E42' = '  +0.31%  '
    'B43' = 'Bittensor'
    'C43' = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    'D43' = '''275.77'
    'E43' = '  -4.02%  '
    'B44' = 'EnergySwap'
    'C44' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D44' = '''19.89'
    'E44' = '  -0.85%  '
    'D45' = '''0.0972'
    'E45' = '  -2.18%  '
    'D46' = '''0.0536'
    'E46' = '  -1.14%  '
    'E47' = '  +2.04%  '
    'B48' = 'RenderToken'
    'C48' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D48' = '''4.78'
    'E48' = '  -2.92%  '
    'B49' = 'Maker'
    'C49' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D49' = '2.029.28'
    'E49' = '  -5.55%  '
    'E50' = '  -2.72%  '
    'D51' = '''18.96'
    'E51' = '  -1.78%  '
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

Write-Host "Applied $($updates.Count) cell updates"
